# Update column F (dSF) values on the active sheet per repull/push data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 3
    4  = 3
    5  = -1
    6  = 2
    7  = 1
    8  = -2
    9  = -3
    10 = -7
    11 = 9
    12 = 4
    14 = 1
    15 = -2
    16 = -4
    17 = -1
    18 = 1
    19 = 1
    20 = 4
    21 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
